# CryCompanywiseStockReport_1 - stock-take correction pass
# Adjusts sold/issued quantities (col F) and recomputed stock values (col G)
# for the affected line items, then rolls the corrections up through the
# per-company "Sub Total:" rows (col B) and the workbook Grand Total (B930/B931).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Line item corrections: Qty (F) and Value (G) ---
$ws.Cells.Item(17, 6).Value = 47
$ws.Cells.Item(17, 7).Value = 7411.9
$ws.Cells.Item(130, 6).Value = 92
$ws.Cells.Item(130, 7).Value = 4552.16
$ws.Cells.Item(155, 2).Value = 53925   # item code
$ws.Cells.Item(155, 6).Value = 1
$ws.Cells.Item(155, 7).Value = 66.44
$ws.Cells.Item(156, 2).Value = 57756   # item code
$ws.Cells.Item(156, 6).Value = 60
$ws.Cells.Item(156, 7).Value = 3986.4
$ws.Cells.Item(167, 6).Value = 6
$ws.Cells.Item(167, 7).Value = 4876.98
$ws.Cells.Item(181, 6).Value = 35
$ws.Cells.Item(181, 7).Value = 3808.35
$ws.Cells.Item(192, 6).Value = 25
$ws.Cells.Item(192, 7).Value = 1886.5
$ws.Cells.Item(228, 6).Value = 43
$ws.Cells.Item(228, 7).Value = 3483
$ws.Cells.Item(229, 6).Value = 12
$ws.Cells.Item(229, 7).Value = 1230.6
$ws.Cells.Item(252, 6).Value = 205
$ws.Cells.Item(252, 7).Value = 3792.5
$ws.Cells.Item(292, 6).Value = 98
$ws.Cells.Item(292, 7).Value = 13282.92
$ws.Cells.Item(302, 6).Value = 101
$ws.Cells.Item(302, 7).Value = 4671.25
$ws.Cells.Item(306, 6).Value = 90
$ws.Cells.Item(306, 7).Value = 2886.3
$ws.Cells.Item(312, 6).Value = 216
$ws.Cells.Item(312, 7).Value = 7192.8
$ws.Cells.Item(328, 6).Value = 0
$ws.Cells.Item(328, 7).Value = 0
$ws.Cells.Item(330, 6).Value = 15
$ws.Cells.Item(330, 7).Value = 1315.5
$ws.Cells.Item(334, 6).Value = 10
$ws.Cells.Item(334, 7).Value = 1730.6
$ws.Cells.Item(337, 6).Value = 26
$ws.Cells.Item(337, 7).Value = 4989.4
$ws.Cells.Item(338, 6).Value = 5
$ws.Cells.Item(338, 7).Value = 414.7
$ws.Cells.Item(341, 6).Value = 237
$ws.Cells.Item(341, 7).Value = 40951.23
$ws.Cells.Item(342, 6).Value = 6
$ws.Cells.Item(342, 7).Value = 1036.74
$ws.Cells.Item(345, 6).Value = 64
$ws.Cells.Item(345, 7).Value = 5258.24
$ws.Cells.Item(346, 6).Value = 206
$ws.Cells.Item(346, 7).Value = 33145.4
$ws.Cells.Item(350, 6).Value = 245
$ws.Cells.Item(350, 7).Value = 18058.95
$ws.Cells.Item(351, 6).Value = 27
$ws.Cells.Item(351, 7).Value = 4284.09
$ws.Cells.Item(352, 6).Value = 98
$ws.Cells.Item(352, 7).Value = 6808.06
$ws.Cells.Item(353, 6).Value = 210
$ws.Cells.Item(353, 7).Value = 8116.5
$ws.Cells.Item(355, 6).Value = 34
$ws.Cells.Item(355, 7).Value = 4455.36
$ws.Cells.Item(360, 6).Value = 99
$ws.Cells.Item(360, 7).Value = 14210.46
$ws.Cells.Item(361, 6).Value = 12
$ws.Cells.Item(361, 7).Value = 1354.08
$ws.Cells.Item(362, 6).Value = 62
$ws.Cells.Item(362, 7).Value = 4090.76
$ws.Cells.Item(363, 6).Value = 1
$ws.Cells.Item(363, 7).Value = 134.78
$ws.Cells.Item(365, 6).Value = 44
$ws.Cells.Item(365, 7).Value = 3655.52
$ws.Cells.Item(369, 6).Value = 30
$ws.Cells.Item(369, 7).Value = 406.8
$ws.Cells.Item(379, 6).Value = 335
$ws.Cells.Item(379, 7).Value = 7745.2
$ws.Cells.Item(384, 6).Value = 9
$ws.Cells.Item(384, 7).Value = 2057.94
$ws.Cells.Item(386, 6).Value = 22
$ws.Cells.Item(386, 7).Value = 2733.94
$ws.Cells.Item(387, 6).Value = 22
$ws.Cells.Item(387, 7).Value = 1393.7
$ws.Cells.Item(392, 6).Value = 64
$ws.Cells.Item(392, 7).Value = 5486.08
$ws.Cells.Item(399, 6).Value = 37
$ws.Cells.Item(399, 7).Value = 4603.54
$ws.Cells.Item(401, 6).Value = 830
$ws.Cells.Item(401, 7).Value = 48721
$ws.Cells.Item(404, 6).Value = 40
$ws.Cells.Item(404, 7).Value = 8665.2
$ws.Cells.Item(415, 6).Value = 62
$ws.Cells.Item(415, 7).Value = 5152.2
$ws.Cells.Item(417, 6).Value = 650
$ws.Cells.Item(417, 7).Value = 111364.5
$ws.Cells.Item(418, 6).Value = 235
$ws.Cells.Item(418, 7).Value = 35524.95
$ws.Cells.Item(419, 6).Value = 5
$ws.Cells.Item(419, 7).Value = 2092.45
$ws.Cells.Item(421, 6).Value = 20
$ws.Cells.Item(421, 7).Value = 3206.4
$ws.Cells.Item(423, 6).Value = 32
$ws.Cells.Item(423, 7).Value = 1972.48
$ws.Cells.Item(429, 6).Value = 448
$ws.Cells.Item(429, 7).Value = 26642.56
$ws.Cells.Item(430, 6).Value = 448
$ws.Cells.Item(430, 7).Value = 18448.64
$ws.Cells.Item(431, 6).Value = 222
$ws.Cells.Item(431, 7).Value = 20563.86
$ws.Cells.Item(432, 6).Value = 501
$ws.Cells.Item(432, 7).Value = 19914.75
$ws.Cells.Item(434, 6).Value = 140
$ws.Cells.Item(434, 7).Value = 20112.4
$ws.Cells.Item(437, 6).Value = 126
$ws.Cells.Item(437, 7).Value = 23131.08
$ws.Cells.Item(438, 6).Value = 129
$ws.Cells.Item(438, 7).Value = 23681.82
$ws.Cells.Item(440, 6).Value = 40
$ws.Cells.Item(440, 7).Value = 6136.4
$ws.Cells.Item(442, 6).Value = 41
$ws.Cells.Item(442, 7).Value = 6630.93
$ws.Cells.Item(445, 6).Value = 19
$ws.Cells.Item(445, 7).Value = 3666.24
$ws.Cells.Item(471, 6).Value = 13
$ws.Cells.Item(471, 7).Value = 1259.7
$ws.Cells.Item(484, 6).Value = 167
$ws.Cells.Item(484, 7).Value = 4063.11
$ws.Cells.Item(494, 6).Value = 104
$ws.Cells.Item(494, 7).Value = 683.28
$ws.Cells.Item(504, 6).Value = 0
$ws.Cells.Item(504, 7).Value = 0
$ws.Cells.Item(556, 6).Value = 46
$ws.Cells.Item(556, 7).Value = 2581.98
$ws.Cells.Item(558, 6).Value = 10
$ws.Cells.Item(558, 7).Value = 1842.3
$ws.Cells.Item(560, 6).Value = 53
$ws.Cells.Item(560, 7).Value = 2313.45
$ws.Cells.Item(562, 6).Value = 7
$ws.Cells.Item(562, 7).Value = 331.8
$ws.Cells.Item(563, 6).Value = 36
$ws.Cells.Item(563, 7).Value = 3316.32
$ws.Cells.Item(564, 6).Value = 48
$ws.Cells.Item(564, 7).Value = 4937.76
$ws.Cells.Item(567, 6).Value = 10
$ws.Cells.Item(567, 7).Value = 1864.4
$ws.Cells.Item(568, 6).Value = 46
$ws.Cells.Item(568, 7).Value = 3310.62
$ws.Cells.Item(571, 6).Value = 39
$ws.Cells.Item(571, 7).Value = 1122.42
$ws.Cells.Item(574, 6).Value = 9
$ws.Cells.Item(574, 7).Value = 341.82
$ws.Cells.Item(576, 6).Value = 10
$ws.Cells.Item(576, 7).Value = 358.9
$ws.Cells.Item(611, 6).Value = 3
$ws.Cells.Item(611, 7).Value = 9463.23
$ws.Cells.Item(677, 6).Value = 25
$ws.Cells.Item(677, 7).Value = 1336
$ws.Cells.Item(721, 6).Value = 10
$ws.Cells.Item(721, 7).Value = 1038.7
$ws.Cells.Item(724, 6).Value = 98
$ws.Cells.Item(724, 7).Value = 4922.54
$ws.Cells.Item(798, 6).Value = 25
$ws.Cells.Item(798, 7).Value = 5313
$ws.Cells.Item(800, 6).Value = 151
$ws.Cells.Item(800, 7).Value = 12910.5
$ws.Cells.Item(803, 6).Value = 61
$ws.Cells.Item(803, 7).Value = 4365.16
$ws.Cells.Item(805, 6).Value = 12
$ws.Cells.Item(805, 7).Value = 964.56
$ws.Cells.Item(806, 6).Value = 104
$ws.Cells.Item(806, 7).Value = 7442.24
$ws.Cells.Item(843, 6).Value = 271
$ws.Cells.Item(843, 7).Value = 22102.76
$ws.Cells.Item(844, 6).Value = 225
$ws.Cells.Item(844, 7).Value = 10768.5
$ws.Cells.Item(845, 6).Value = 12
$ws.Cells.Item(845, 7).Value = 978.72
$ws.Cells.Item(846, 6).Value = 91
$ws.Cells.Item(846, 7).Value = 14061.32
$ws.Cells.Item(847, 6).Value = 191
$ws.Cells.Item(847, 7).Value = 15577.96
$ws.Cells.Item(848, 6).Value = 364
$ws.Cells.Item(848, 7).Value = 48448.4
$ws.Cells.Item(852, 6).Value = 127
$ws.Cells.Item(852, 7).Value = 2758.44
$ws.Cells.Item(853, 6).Value = 188
$ws.Cells.Item(853, 7).Value = 7008.64
$ws.Cells.Item(855, 6).Value = 384
$ws.Cells.Item(855, 7).Value = 26711.04
$ws.Cells.Item(863, 6).Value = 434
$ws.Cells.Item(863, 7).Value = 62496
$ws.Cells.Item(865, 6).Value = 305
$ws.Cells.Item(865, 7).Value = 36816.55
$ws.Cells.Item(893, 6).Value = 9
$ws.Cells.Item(893, 7).Value = 47444.67
$ws.Cells.Item(899, 6).Value = 8
$ws.Cells.Item(899, 7).Value = 113249.28
$ws.Cells.Item(912, 6).Value = 1810
$ws.Cells.Item(912, 7).Value = 295229.1

# --- Per-company "Sub Total:" rollups (col B) ---
$ws.Cells.Item(19, 2).Value = 55816.69
$ws.Cells.Item(140, 2).Value = 58838.57
$ws.Cells.Item(188, 2).Value = 208619.6
$ws.Cells.Item(205, 2).Value = 7483.91
$ws.Cells.Item(230, 2).Value = 19313.11
$ws.Cells.Item(259, 2).Value = 7642.67
$ws.Cells.Item(331, 2).Value = 213887.17
$ws.Cells.Item(435, 2).Value = 699128.49
$ws.Cells.Item(453, 2).Value = 110695.53
$ws.Cells.Item(477, 2).Value = 34937.53
$ws.Cells.Item(507, 2).Value = 124600.01
$ws.Cells.Item(577, 2).Value = 38610.83
$ws.Cells.Item(617, 2).Value = 243653.48
$ws.Cells.Item(682, 2).Value = 1737.39
$ws.Cells.Item(735, 2).Value = 23302.95
$ws.Cells.Item(807, 2).Value = 57776.48
$ws.Cells.Item(867, 2).Value = 473210.04
$ws.Cells.Item(902, 2).Value = 452996.94
$ws.Cells.Item(918, 2).Value = 328971.74

# --- Workbook-wide Grand Total (col B) ---
$ws.Cells.Item(930, 2).Value = 6011930.83
$ws.Cells.Item(931, 2).Value = 6011930.83
